$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 65
$ws.Range("F4").Value = 1590
$ws.Range("F6").Value = 619
$ws.Range("F7").Value = 1103
$ws.Range("F8").Value = 1547
$ws.Range("F12").Value = 1459
$ws.Range("F13").Value = 3081
$ws.Range("F14").Value = 618
$ws.Range("F16").Value = 1803
$ws.Range("F17").Value = 850
$ws.Range("F18").Value = 277
$ws.Range("F23").Value = 7
$ws.Range("F24").Value = 1209
$ws.Range("F25").Value = 403
$ws.Range("F26").Value = 455
$ws.Range("F27").Value = 110
$ws.Range("F28").Value = 4784
$ws.Range("F29").Value = 35
$ws.Range("F31").Value = 572
$ws.Range("F32").Value = 1652
$ws.Range("F34").Value = 135

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 56
$ws.Range("F7").Value = 71

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 43

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 65
$ws.Range("F3").Value = 43
$ws.Range("F4").Value = 56
$ws.Range("F10").Value = 71
$ws.Range("F12").Value = 1590
$ws.Range("F14").Value = 619
$ws.Range("F15").Value = 1103
$ws.Range("F16").Value = 1547
$ws.Range("F21").Value = 1459
$ws.Range("F22").Value = 3081
$ws.Range("F23").Value = 618
$ws.Range("F25").Value = 1803
$ws.Range("F26").Value = 851
$ws.Range("F27").Value = 277
$ws.Range("F33").Value = 7
$ws.Range("F35").Value = 1209
$ws.Range("F36").Value = 403
$ws.Range("F37").Value = 455
$ws.Range("F38").Value = 110
$ws.Range("F39").Value = 4784
$ws.Range("F40").Value = 35
$ws.Range("F42").Value = 572
$ws.Range("F43").Value = 1652
$ws.Range("F47").Value = 135
